# Rebuild the Sheet1 backlog table per the target layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate so removed rows/cells truly disappear.
$ws.Rows("1:20").Delete()

# ---- Header row ----
$ws.Range("A1").Value = "Feature"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Techncial Features"
$ws.Range("D1").Value = "Status"
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("B1:C1").WrapText = $true

# ---- Row 2 ----
$ws.Range("A2").Value = "Trailer stop loss"
$ws.Range("B2").Value = "Trailer stop loss for buy and sell orders, and with stop loss in the other direction for binance"
$ws.Range("C2").Value = "service"
$ws.Range("D2").Value = "Completed"
$ws.Range("B2:C2").WrapText = $true
$ws.Rows("2:2").RowHeight = 29

# ---- Row 3 ----
$ws.Range("C3").Value = "saving log of trades using pandas"
$ws.Range("D3").Value = "Completed"
$ws.Range("C3").WrapText = $true

# (row 4 intentionally left blank)

# ---- Row 5 ----
$ws.Range("A5").Value = "MA & EMA cross bot"
$ws.Range("B5").Value = "Automate buy and sell based on the cross of the moving averages"
$ws.Range("C5").Value = "Basic bot"
$ws.Range("D5").Value = "In progress"
$ws.Range("B5:C5").WrapText = $true

# ---- Row 6 ----
$ws.Range("C6").Value = "History emulation"
$ws.Range("D6").Value = "Discovery"
$ws.Range("C6").WrapText = $true

# ---- Row 7 ----
$ws.Range("C7").Value = "Demo run feature"
$ws.Range("D7").Value = "Discovery"
$ws.Range("C7").WrapText = $true

# ---- Row 8 ----
$ws.Range("C8").Value = "run with actual sell and buy on market price"
$ws.Range("D8").Value = "Discovery"
$ws.Range("C8").WrapText = $true
$ws.Rows("8:8").RowHeight = 29

# ---- Row 9 ----
$ws.Range("A9").Value = "Momentum Indicators bot"
$ws.Range("B9").Value = "Automate buy and sell based on indicators like the RSA or MACD."
$ws.Range("D9").Value = "Discovery"
$ws.Range("B9").WrapText = $true

# ---- Row 10 ----
$ws.Range("A10").Value = "Mix Bot"
$ws.Range("B10").Value = "Combine overlap bot with momentum bot"
$ws.Range("D10").Value = "Discovery"
$ws.Range("B10").WrapText = $true

# ---- Row 11 ----
$ws.Range("A11").Value = "ICT bot"
$ws.Range("B11").Value = "Candle pattern recognition for ICT order blocks."
$ws.Range("D11").Value = "Discovery"
$ws.Range("B11").WrapText = $true

# ---- Row 12 ----
$ws.Range("A12").Value = "Ice berg orders bot"
$ws.Range("B12").Value = "Ice berg bot"
$ws.Range("D12").Value = "Discovery"
$ws.Range("B12").WrapText = $true

# ---- Row 13 ----
$ws.Range("A13").Value = "Arbitrage bot"
$ws.Range("D13").Value = "Discovery"

# ---- Row 14 ----
$ws.Range("A14").Value = "Build GUI"
$ws.Range("D14").Value = "Discovery"

# ---- Row 15 ----
$ws.Range("A15").Value = "Statistics and reports"
$ws.Range("D15").Value = "Discovery"

# Selection as recorded in the saved workbook.
$ws.Range("D7:D8").Select()
